# Weekly price update: insert a new daily price record for "Ají" (Inferno /
# Primera) at row 39, pushing the existing historical rows (old 39..136)
# down by one (new 40..137). The new row is seeded as a copy of the row
# that lands right below it (so every column besides the price/date fields
# keeps its usual value), then the date, volume, min/max/avg price and
# $/kg fields are overwritten with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 39:136 down to 40:137, opening up a blank row 39.
$ws.Rows(39).Insert()

# Seed the new row 39 with the same shape/content as the row now below it
# (old row 39, now at row 40) so every non-edited column already matches.
$ws.Rows(40).Copy()
$ws.Rows(39).PasteSpecial()

# Overwrite this week's actual reported values.
$ws.Cells.Item(39, 4).Value = 45044   # D39 Fecha
$ws.Cells.Item(39, 10).Value = 280    # J39 Volumen
$ws.Cells.Item(39, 11).Value = 12000  # K39 Precio minimo
$ws.Cells.Item(39, 12).Value = 13000  # L39 Precio maximo
$ws.Cells.Item(39, 13).Value = 12643  # M39 Precio promedio ponderado
$ws.Cells.Item(39, 16).Value = 843    # P39 Precio $/Kg
